# Update Name of Algo
# Applies the updated RandomForest imputed values for the terrestrial_mammals
# combination_2_ABCDE / AD / 20 / seed2 result data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = -7.255500000000003
$ws.Range("A3").Value = -22.31000000000001
$ws.Range("A14").Value = -21.7461
$ws.Range("A16").Value = -21.66069999999998
$ws.Range("D18").Value = -8.872900000000001
$ws.Range("A21").Value = -20.31539999999998
$ws.Range("A23").Value = -20.56499999999998
$ws.Range("D24").Value = -7.3515
$ws.Range("A25").Value = -21.75079999999999
$ws.Range("D25").Value = -8.008099999999997
$ws.Range("A26").Value = -21.31079999999997
$ws.Range("D27").Value = -8.988199999999999
$ws.Range("A29").Value = -20.96479999999998
$ws.Range("D30").Value = -7.397500000000007
$ws.Range("D31").Value = -8.304400000000003
$ws.Range("D39").Value = -8.484799999999998
$ws.Range("A40").Value = -20.29489999999999
$ws.Range("D42").Value = -8.917899999999999
$ws.Range("D48").Value = -7.334099999999999
$ws.Range("D51").Value = -7.900299999999999
$ws.Range("D52").Value = -7.6454
$ws.Range("A53").Value = -21.7968
$ws.Range("D55").Value = -9.015300000000002
$ws.Range("D56").Value = -7.771999999999996
$ws.Range("A57").Value = -22.67910000000001
$ws.Range("D57").Value = -8.638100000000001
$ws.Range("A59").Value = -22.7455
$ws.Range("D60").Value = -7.806299999999993
$ws.Range("A65").Value = -21.77649999999998
$ws.Range("A69").Value = -21.60250000000001
$ws.Range("D73").Value = -7.259999999999996
$ws.Range("D74").Value = -8.686900000000003
$ws.Range("A79").Value = -20.42340000000002
$ws.Range("A83").Value = -21.81109999999999
$ws.Range("D89").Value = -5.8355
$ws.Range("D90").Value = -7.897200000000003
$ws.Range("A91").Value = -21.42710000000002
$ws.Range("D92").Value = -5.974999999999998
$ws.Range("A93").Value = -21.05649999999998
$ws.Range("A100").Value = -21.8385
